$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.899.84'
$ws.Range('E2').Value = '  -0.25%  '
$ws.Range('D3').Value = '1.898.00'
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = '0.7918'
$ws.Range('E5').Value = '  -4.48%  '
$ws.Range('D6').Value = '244.06'
$ws.Range('E6').Value = '  +0.97%  '
$ws.Range('D7').Value = '1.001'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.3166'
$ws.Range('E8').Value = '  -3.29%  '
$ws.Range('D9').Value = '25.41'
$ws.Range('E9').Value = '  -4.03%  '
$ws.Range('D10').Value = '0.07165'
$ws.Range('E10').Value = '  +1.93%  '
$ws.Range('D11').Value = '0.08112'
$ws.Range('D12').Value = '5.591'
$ws.Range('E12').Value = '  +6.75%  '
$ws.Range('D13').Value = '0.7682'
$ws.Range('E13').Value = '  +1.14%  '
$ws.Range('D14').Value = '1.923.09'
$ws.Range('E14').Value = '  +1.34%  '
$ws.Range('D15').Value = '92.63'
$ws.Range('E15').Value = '  +0.59%  '
$ws.Range('D16').Value = '6.173'
$ws.Range('E16').Value = '  +5.40%  '
$ws.Range('D17').Value = '29.917.62'
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '13.95'
$ws.Range('E18').Value = '  -1.00%  '
$ws.Range('D19').Value = '244.59'
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('D20').Value = '0.000007784'
$ws.Range('E20').Value = '  +0.56%  '
$ws.Range('D21').Value = '8.298'
$ws.Range('E21').Value = '  +19.39%  '
$ws.Range('D22').Value = '2.170.14'
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('D23').Value = '1.002'
$ws.Range('E23').Value = '  +0.12%  '
$ws.Range('D24').Value = '1.002'
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1670'
$ws.Range('E25').Value = '  -3.72%  '
$ws.Range('D26').Value = '9.491'
$ws.Range('E26').Value = '  +2.70%  '
$ws.Range('D27').Value = '163.95'
$ws.Range('E27').Value = '  -0.99%  '
$ws.Range('D28').Value = '18.75'
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('D29').Value = '2.067'
$ws.Range('E29').Value = '  -1.01%  '
$ws.Range('D30').Value = '1.401'
$ws.Range('E30').Value = '  +2.97%  '
$ws.Range('E31').Value = '  +2.28%  '
$ws.Range('D32').Value = '4.496'
$ws.Range('E32').Value = '  +5.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05610'
$ws.Range('E33').Value = '  -5.62%  '
$ws.Range('D34').Value = '4.098'
$ws.Range('E34').Value = '  +0.85%  '
$ws.Range('D35').Value = '1.279'
$ws.Range('E35').Value = '  +1.04%  '
$ws.Range('D36').Value = '0.7427'
$ws.Range('E36').Value = '  +1.73%  '
$ws.Range('D38').Value = '2.635'
$ws.Range('E38').Value = '  -3.19%  '
$ws.Range('D39').Value = '0.01935'
$ws.Range('E39').Value = '  +1.14%  '
$ws.Range('D40').Value = '2.784'
$ws.Range('E40').Value = '  +0.22%  '
$ws.Range('D41').Value = '1.172.44'
$ws.Range('E41').Value = '  +19.07%  '
$ws.Range('D42').Value = '74.65'
$ws.Range('E42').Value = '  +3.29%  '
$ws.Range('D43').Value = '0.4429'
$ws.Range('E43').Value = '  -0.05%  '
$ws.Range('E44').Value = '  +2.08%  '
$ws.Range('D45').Value = '0.8534'
$ws.Range('E45').Value = '  +0.31%  '
$ws.Range('D46').Value = '104.79'
$ws.Range('E46').Value = '  +2.84%  '
$ws.Range('D47').Value = '1.001'
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('E48').Value = '  -0.40%  '
$ws.Range('D49').Value = '10.02'
$ws.Range('E49').Value = '  +2.03%  '
$ws.Range('D50').Value = '7.473'
$ws.Range('E50').Value = '  -0.85%  '
$ws.Range('E51').Value = '  +10.41%  '
